# Generate Report for Archive
#
# A new handoff round completed for f308f5f0-ec79-4bf6-a5d5-18ff80e7476b
# (status moved from "Ready for handoff" to "In Translation", with a new
# handoff file/timestamp), so the generated status report now lists that
# file in row 6 (ahead of 06735111-ee76-4e32-890b-cd65bbea6fe1, which
# drops to row 7 but keeps its prior "Ready for handoff" status).
#
# The row's hyperlink *targets* (r:id / the underlying URL) stay anchored
# to the row position - only the cell values and the hyperlinks' visible
# text move with the data.

$wb = $excel.ActiveWorkbook

function Set-RowSixSeven($SheetName, $Row6A, $Row6B, $Row6C, $Row6D, $Row7A, $Row7B, $Row7C, $Row7D, $HasD) {
    $ws = $wb.Worksheets.Item($SheetName)

    # --- cell values -----------------------------------------------------
    $ws.Range("A6").Value = $Row6A
    $ws.Range("B6").Value = $Row6B
    $ws.Range("C6").Value = $Row6C
    if ($HasD) { $ws.Range("D6").Value = $Row6D }

    $ws.Range("A7").Value = $Row7A
    $ws.Range("B7").Value = $Row7B
    $ws.Range("C7").Value = $Row7C
    if ($HasD) { $ws.Range("D7").Value = $Row7D }

    # --- hyperlink display text (targets stay put on the row) -----------
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq '$A$6') { $h.TextToDisplay = $Row6A }
        elseif ($addr -eq '$A$7') { $h.TextToDisplay = $Row7A }
        elseif ($addr -eq '$C$6') { $h.TextToDisplay = $Row6C }
        elseif ($addr -eq '$C$7') { $h.TextToDisplay = $Row7C }
    }
}

# Overview sheet: columns A (file), B (zh-cn status), C (de-de status)
Set-RowSixSeven "Overview" `
    "f308f5f0-ec79-4bf6-a5d5-18ff80e7476b.md" "In Translation" "In Translation" $null `
    "06735111-ee76-4e32-890b-cd65bbea6fe1.md" "Ready for handoff" "Ready for handoff" $null `
    $false

# zh-cn sheet: A (file), B (status), C (handoff file), D (handoff datetime)
Set-RowSixSeven "zh-cn" `
    "f308f5f0-ec79-4bf6-a5d5-18ff80e7476b.md" "In Translation" `
    "f308f5f0-ec79-4bf6-a5d5-18ff80e7476b.b8f21af0909c18c0eb6be2afb8bd10a4e710c4a4.zh-cn.xlf" "2016-01-28 05:13:29" `
    "06735111-ee76-4e32-890b-cd65bbea6fe1.md" "Ready for handoff" `
    "06735111-ee76-4e32-890b-cd65bbea6fe1.0733acfe6952d111aa43861c80d3c0912bc917f6.zh-cn.xlf" "2016-01-28 05:09:14" `
    $true

# de-de sheet: A (file), B (status), C (handoff file), D (handoff datetime)
Set-RowSixSeven "de-de" `
    "f308f5f0-ec79-4bf6-a5d5-18ff80e7476b.md" "In Translation" `
    "f308f5f0-ec79-4bf6-a5d5-18ff80e7476b.b8f21af0909c18c0eb6be2afb8bd10a4e710c4a4.de-de.xlf" "2016-01-28 05:13:40" `
    "06735111-ee76-4e32-890b-cd65bbea6fe1.md" "Ready for handoff" `
    "06735111-ee76-4e32-890b-cd65bbea6fe1.0733acfe6952d111aa43861c80d3c0912bc917f6.de-de.xlf" "2016-01-28 05:09:24" `
    $true
